# Update 1axis tracking solar persistence
# - Bump the model name/version in "Main Parameters"
# - Update site Latitude/Longitude in "Main Parameters"
# - Point the Timeseries Attributes input files at the new PSE-scaled CSVs

$wb = $excel.ActiveWorkbook

# --- Timeseries Attributes --------------------------------------------
$tsAttrs = $wb.Worksheets.Item("Timeseries Attributes")
$tsAttrs.Range("B2").Value = "PSE_RTD_load_forecast_scaled.csv"
$tsAttrs.Range("B3").Value = "PSE_RTPD_load_forecast_scaled.csv"
$tsAttrs.Range("B4").Value = "PSE_solar_5_minute_actuals.csv"
$tsAttrs.Range("B5").Value = "PSE_wind_5_minute_actuals.csv"

# --- Main Parameters -------------------------------------------------
$mainParams = $wb.Worksheets.Item("Main Parameters")
$mainParams.Range("B2").Value = "RESERVE_PSE_6.0"
$mainParams.Range("B3").Value = 45.9
$mainParams.Range("B4").Value = -106.62
